# Update header info on all three sheets (new scrape timestamp)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Range("A2").Value = "Última actualización: 11:47:13"
$ws1.Range("A3").Value = "Total filas: 204"
$ws2.Range("A2").Value = "Última actualización: 11:47:13"
$ws3.Range("A2").Value = "Última actualización: 11:47:13"

# Refresh scraped schedule rows on sheet "LP1912" (rows 81-209)
$ws1.Cells.Item(81, 1).Value = "07:12:53"
$ws1.Cells.Item(81, 2).Value = "09:02"
$ws1.Cells.Item(81, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(81, 4).Value = 110
$ws1.Cells.Item(81, 5).Value = "LP1912"
$ws1.Cells.Item(82, 1).Value = "07:36:59"
$ws1.Cells.Item(82, 2).Value = "09:02"
$ws1.Cells.Item(82, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(82, 4).Value = 86
$ws1.Cells.Item(82, 5).Value = "LP1912"
$ws1.Cells.Item(107, 1).Value = "08:11:27"
$ws1.Cells.Item(107, 2).Value = "10:03"
$ws1.Cells.Item(107, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(107, 4).Value = 112
$ws1.Cells.Item(107, 5).Value = "LP1912"
$ws1.Cells.Item(108, 1).Value = "09:21:49"
$ws1.Cells.Item(108, 2).Value = "10:03"
$ws1.Cells.Item(108, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(108, 4).Value = 42
$ws1.Cells.Item(108, 5).Value = "LP1912"
$ws1.Cells.Item(158, 1).Value = "10:04:17"
$ws1.Cells.Item(158, 2).Value = "11:25"
$ws1.Cells.Item(158, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(158, 4).Value = 81
$ws1.Cells.Item(158, 5).Value = "LP1912"
$ws1.Cells.Item(159, 1).Value = "10:36:18"
$ws1.Cells.Item(159, 2).Value = "11:25"
$ws1.Cells.Item(159, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(159, 4).Value = 49
$ws1.Cells.Item(159, 5).Value = "LP1912"
$ws1.Cells.Item(175, 1).Value = "11:47:13"
$ws1.Cells.Item(175, 2).Value = "11:57"
$ws1.Cells.Item(175, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(175, 4).Value = 10
$ws1.Cells.Item(175, 5).Value = "LP1912"
$ws1.Cells.Item(176, 1).Value = "11:34:25"
$ws1.Cells.Item(176, 2).Value = "12:03"
$ws1.Cells.Item(176, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(176, 4).Value = 29
$ws1.Cells.Item(176, 5).Value = "LP1912"
$ws1.Cells.Item(177, 1).Value = "11:11:31"
$ws1.Cells.Item(177, 2).Value = "12:05"
$ws1.Cells.Item(177, 3).Value = "17_ROMERO"
$ws1.Cells.Item(177, 4).Value = 54
$ws1.Cells.Item(177, 5).Value = "LP1912"
$ws1.Cells.Item(178, 1).Value = "11:11:31"
$ws1.Cells.Item(178, 2).Value = "12:06"
$ws1.Cells.Item(178, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(178, 4).Value = 55
$ws1.Cells.Item(178, 5).Value = "LP1912"
$ws1.Cells.Item(179, 1).Value = "10:48:14"
$ws1.Cells.Item(179, 2).Value = "12:07"
$ws1.Cells.Item(179, 3).Value = "14_ABASTO"
$ws1.Cells.Item(179, 4).Value = 79
$ws1.Cells.Item(179, 5).Value = "LP1912"
$ws1.Cells.Item(180, 1).Value = "11:47:13"
$ws1.Cells.Item(180, 2).Value = "12:07"
$ws1.Cells.Item(180, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(180, 4).Value = 20
$ws1.Cells.Item(180, 5).Value = "LP1912"
$ws1.Cells.Item(181, 1).Value = "11:11:31"
$ws1.Cells.Item(181, 2).Value = "12:17"
$ws1.Cells.Item(181, 3).Value = "15_ABASTO"
$ws1.Cells.Item(181, 4).Value = 66
$ws1.Cells.Item(181, 5).Value = "LP1912"
$ws1.Cells.Item(182, 1).Value = "11:11:31"
$ws1.Cells.Item(182, 2).Value = "12:18"
$ws1.Cells.Item(182, 3).Value = "10_OLMOS"
$ws1.Cells.Item(182, 4).Value = 67
$ws1.Cells.Item(182, 5).Value = "LP1912"
$ws1.Cells.Item(183, 1).Value = "11:34:25"
$ws1.Cells.Item(183, 2).Value = "12:20"
$ws1.Cells.Item(183, 3).Value = "17_ROMERO"
$ws1.Cells.Item(183, 4).Value = 46
$ws1.Cells.Item(183, 5).Value = "LP1912"
$ws1.Cells.Item(184, 1).Value = "11:47:13"
$ws1.Cells.Item(184, 2).Value = "12:21"
$ws1.Cells.Item(184, 3).Value = "17_ROMERO"
$ws1.Cells.Item(184, 4).Value = 34
$ws1.Cells.Item(184, 5).Value = "LP1912"
$ws1.Cells.Item(185, 1).Value = "10:36:18"
$ws1.Cells.Item(185, 2).Value = "12:29"
$ws1.Cells.Item(185, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(185, 4).Value = 113
$ws1.Cells.Item(185, 5).Value = "LP1912"
$ws1.Cells.Item(186, 1).Value = "10:36:18"
$ws1.Cells.Item(186, 2).Value = "12:30"
$ws1.Cells.Item(186, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(186, 4).Value = 114
$ws1.Cells.Item(186, 5).Value = "LP1912"
$ws1.Cells.Item(187, 1).Value = "10:36:18"
$ws1.Cells.Item(187, 2).Value = "12:31"
$ws1.Cells.Item(187, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(187, 4).Value = 115
$ws1.Cells.Item(187, 5).Value = "LP1912"
$ws1.Cells.Item(188, 1).Value = "10:48:14"
$ws1.Cells.Item(188, 2).Value = "12:31"
$ws1.Cells.Item(188, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(188, 4).Value = 103
$ws1.Cells.Item(188, 5).Value = "LP1912"
$ws1.Cells.Item(189, 1).Value = "10:55:25"
$ws1.Cells.Item(189, 2).Value = "12:36"
$ws1.Cells.Item(189, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(189, 4).Value = 101
$ws1.Cells.Item(189, 5).Value = "LP1912"
$ws1.Cells.Item(190, 1).Value = "10:48:14"
$ws1.Cells.Item(190, 2).Value = "12:37"
$ws1.Cells.Item(190, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(190, 4).Value = 109
$ws1.Cells.Item(190, 5).Value = "LP1912"
$ws1.Cells.Item(191, 1).Value = "10:48:14"
$ws1.Cells.Item(191, 2).Value = "12:40"
$ws1.Cells.Item(191, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(191, 4).Value = 112
$ws1.Cells.Item(191, 5).Value = "LP1912"
$ws1.Cells.Item(192, 1).Value = "10:55:25"
$ws1.Cells.Item(192, 2).Value = "12:42"
$ws1.Cells.Item(192, 3).Value = "14_ABASTO"
$ws1.Cells.Item(192, 4).Value = 107
$ws1.Cells.Item(192, 5).Value = "LP1912"
$ws1.Cells.Item(193, 1).Value = "10:55:25"
$ws1.Cells.Item(193, 2).Value = "12:43"
$ws1.Cells.Item(193, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(193, 4).Value = 108
$ws1.Cells.Item(193, 5).Value = "LP1912"
$ws1.Cells.Item(194, 1).Value = "10:48:14"
$ws1.Cells.Item(194, 2).Value = "12:43"
$ws1.Cells.Item(194, 3).Value = "14_ABASTO"
$ws1.Cells.Item(194, 4).Value = 115
$ws1.Cells.Item(194, 5).Value = "LP1912"
$ws1.Cells.Item(195, 1).Value = "11:47:13"
$ws1.Cells.Item(195, 2).Value = "12:51"
$ws1.Cells.Item(195, 3).Value = "15_ABASTO"
$ws1.Cells.Item(195, 4).Value = 64
$ws1.Cells.Item(195, 5).Value = "LP1912"
$ws1.Cells.Item(196, 1).Value = "11:11:31"
$ws1.Cells.Item(196, 2).Value = "12:54"
$ws1.Cells.Item(196, 3).Value = "15X38_ABASTO"
$ws1.Cells.Item(196, 4).Value = 103
$ws1.Cells.Item(196, 5).Value = "LP1912"
$ws1.Cells.Item(197, 1).Value = "11:11:31"
$ws1.Cells.Item(197, 2).Value = "13:01"
$ws1.Cells.Item(197, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(197, 4).Value = 110
$ws1.Cells.Item(197, 5).Value = "LP1912"
$ws1.Cells.Item(198, 1).Value = "11:47:13"
$ws1.Cells.Item(198, 2).Value = "13:05"
$ws1.Cells.Item(198, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(198, 4).Value = 78
$ws1.Cells.Item(198, 5).Value = "LP1912"
$ws1.Cells.Item(199, 1).Value = "11:11:31"
$ws1.Cells.Item(199, 2).Value = "13:06"
$ws1.Cells.Item(199, 3).Value = "14_ABASTO"
$ws1.Cells.Item(199, 4).Value = 115
$ws1.Cells.Item(199, 5).Value = "LP1912"
$ws1.Cells.Item(200, 1).Value = "11:47:13"
$ws1.Cells.Item(200, 2).Value = "13:07"
$ws1.Cells.Item(200, 3).Value = "14_ABASTO"
$ws1.Cells.Item(200, 4).Value = 80
$ws1.Cells.Item(200, 5).Value = "LP1912"
$ws1.Cells.Item(201, 1).Value = "11:34:25"
$ws1.Cells.Item(201, 2).Value = "13:11"
$ws1.Cells.Item(201, 3).Value = "215_ALUAR"
$ws1.Cells.Item(201, 4).Value = 97
$ws1.Cells.Item(201, 5).Value = "LP1912"
$ws1.Cells.Item(202, 1).Value = "11:47:13"
$ws1.Cells.Item(202, 2).Value = "13:11"
$ws1.Cells.Item(202, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(202, 4).Value = 84
$ws1.Cells.Item(202, 5).Value = "LP1912"
$ws1.Cells.Item(203, 1).Value = "11:34:25"
$ws1.Cells.Item(203, 2).Value = "13:18"
$ws1.Cells.Item(203, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(203, 4).Value = 104
$ws1.Cells.Item(203, 5).Value = "LP1912"
$ws1.Cells.Item(204, 1).Value = "11:47:13"
$ws1.Cells.Item(204, 2).Value = "13:19"
$ws1.Cells.Item(204, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(204, 4).Value = 92
$ws1.Cells.Item(204, 5).Value = "LP1912"
$ws1.Cells.Item(205, 1).Value = "11:34:25"
$ws1.Cells.Item(205, 2).Value = "13:21"
$ws1.Cells.Item(205, 3).Value = "17_ROMERO"
$ws1.Cells.Item(205, 4).Value = 107
$ws1.Cells.Item(205, 5).Value = "LP1912"
$ws1.Cells.Item(206, 1).Value = "11:34:25"
$ws1.Cells.Item(206, 2).Value = "13:30"
$ws1.Cells.Item(206, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(206, 4).Value = 116
$ws1.Cells.Item(206, 5).Value = "LP1912"
$ws1.Cells.Item(207, 1).Value = "11:47:13"
$ws1.Cells.Item(207, 2).Value = "13:30"
$ws1.Cells.Item(207, 3).Value = "10_OLMOS"
$ws1.Cells.Item(207, 4).Value = 103
$ws1.Cells.Item(207, 5).Value = "LP1912"
$ws1.Cells.Item(208, 1).Value = "11:47:13"
$ws1.Cells.Item(208, 2).Value = "13:31"
$ws1.Cells.Item(208, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(208, 4).Value = 104
$ws1.Cells.Item(208, 5).Value = "LP1912"
$ws1.Cells.Item(209, 1).Value = "11:47:13"
$ws1.Cells.Item(209, 2).Value = "13:40"
$ws1.Cells.Item(209, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(209, 4).Value = 113
$ws1.Cells.Item(209, 5).Value = "LP1912"
